# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Re-sorted "Pakistan" above "Mexico" (row 16/17) following updated totals
#  - Refreshed numeric Covid figures for several countries
#  - Updated the "Datos actualizados" timestamp string (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: refreshed timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 07:08"

# --- Rows 16 & 17: Pakistan overtakes Mexico in total cases -----------------
$ws.Range("A16").Value = "Pakistan"
$ws.Range("B16").Value = 171666
$ws.Range("C16").Value = 6604
$ws.Range("D16").Value = 63504
$ws.Range("E16").Value = 104780
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 153
$ws.Range("H16").Value = 3382

$ws.Range("A17").Value = "Mexico"
$ws.Range("B17").Value = 170485
$ws.Range("C17").Value = 5030
$ws.Range("D17").Value = 127332
$ws.Range("E17").Value = 22759
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 647
$ws.Range("H17").Value = 20394

# --- Row 56: Kazajistan -------------------------------------------------------
$ws.Range("B56").Value = 16779
$ws.Range("C56").Value = 428
$ws.Range("D56").Value = 10411
$ws.Range("E56").Value = 6255

# --- Row 73: Australia --------------------------------------------------------
$ws.Range("B73").Value = 7436
$ws.Range("C73").Value = 27
$ws.Range("E73").Value = 453

# --- Row 80: Haiti -------------------------------------------------------------
$ws.Range("B80").Value = 4980
$ws.Range("C80").Value = 64
$ws.Range("E80").Value = 4869
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 87

# --- Row 94: Tailandia ---------------------------------------------------------
$ws.Range("B94").Value = 3147
$ws.Range("C94").Value = 1
$ws.Range("D94").Value = 3018
$ws.Range("E94").Value = 71

# --- Row 185: Butan --------------------------------------------------------------
$ws.Range("B185").Value = 68
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 28
$ws.Range("E185").Value = 40

# --- Row 208: Islas Turcas y Caicos ----------------------------------------------
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# --- Row 209: Santa Sede -----------------------------------------------------------
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
